$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Save" header in H1, reusing the exact header formatting (style)
# from the adjacent "sum" header cell (G1) via a format-only copy/paste.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the value for the new Save column in row 2
$ws.Range("H2").Value = 1
